$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes: thinh / 10 / k62 / abcd
$ws.Range("A2").Value = "thịnh"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = "k62"
$ws.Range("D2").Value = "abcd"

# Row 3 changes: nam / 11 / k62 / abcde
$ws.Range("A3").Value = "nam"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = "k62"
$ws.Range("D3").Value = "abcde"

# New rows 4-12: various students, all k62, password cxzv
$ws.Range("A4").Value = "nam"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = "k62"
$ws.Range("D4").Value = "cxzv"

$ws.Range("A5").Value = "linh"
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = "k62"
$ws.Range("D5").Value = "cxzv"

$ws.Range("A6").Value = "minh"
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = "k62"
$ws.Range("D6").Value = "cxzv"

$ws.Range("A7").Value = "hiếu"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = "k62"
$ws.Range("D7").Value = "cxzv"

$ws.Range("A8").Value = "văn a"
$ws.Range("B8").Value = 16
$ws.Range("C8").Value = "k62"
$ws.Range("D8").Value = "cxzv"

$ws.Range("A9").Value = "văn b"
$ws.Range("B9").Value = 17
$ws.Range("C9").Value = "k62"
$ws.Range("D9").Value = "cxzv"

$ws.Range("A10").Value = "văn c"
$ws.Range("B10").Value = 18
$ws.Range("C10").Value = "k62"
$ws.Range("D10").Value = "cxzv"

$ws.Range("A11").Value = "dương"
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = "k62"
$ws.Range("D11").Value = "cxzv"

$ws.Range("A12").Value = "dũng"
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = "k62"
$ws.Range("D12").Value = "cxzv"

# New admin row 13: admin / 1 / none / admin
$ws.Range("A13").Value = "admin"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "none"
$ws.Range("D13").Value = "admin"

# Match the cell selection left by the author on save
[void]$ws.Range("D13").Select()
